# Refresh cryptos list (price + 1h volume change) per the GitHub Actions data pull.
# Leading apostrophe forces text so dot-separated numbers (e.g. "27.964.60")
# and trailing zeros (e.g. "8.690") are preserved exactly, matching the source feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.964.60"
$ws.Range("E2").Value = "'  +4.73%  "

$ws.Range("D3").Value = "'1.781.15"
$ws.Range("E3").Value = "'  +3.29%  "

$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "'  +0.24%  "

$ws.Range("D5").Value = "'244.12"
$ws.Range("E5").Value = "'  +1.00%  "

$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "'  +0.26%  "

$ws.Range("D7").Value = "'0.4895"
$ws.Range("E7").Value = "'  -0.73%  "

$ws.Range("D8").Value = "'0.2673"
$ws.Range("E8").Value = "'  +2.39%  "

$ws.Range("D9").Value = "'0.06256"
$ws.Range("E9").Value = "'  +0.52%  "

$ws.Range("D10").Value = "'1.778.38"
$ws.Range("E10").Value = "'  +3.14%  "

$ws.Range("D11").Value = "'16.34"
$ws.Range("E11").Value = "'  +3.16%  "

$ws.Range("D12").Value = "'0.07031"
$ws.Range("E12").Value = "'  +0.54%  "

$ws.Range("D13").Value = "'0.6276"
$ws.Range("E13").Value = "'  +2.96%  "

$ws.Range("D14").Value = "'4.625"
$ws.Range("E14").Value = "'  +2.76%  "

$ws.Range("D15").Value = "'79.95"
$ws.Range("E15").Value = "'  +3.70%  "

$ws.Range("B16").Value = "Dai"
$ws.Range("C16").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "'  +0.30%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'27.922.18"
$ws.Range("E17").Value = "'  +5.35%  "

$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "'  +0.33%  "

$ws.Range("D19").Value = "'0.000007221"
$ws.Range("E19").Value = "'  +0.36%  "

$ws.Range("D20").Value = "'11.94"
$ws.Range("E20").Value = "'  +4.70%  "

$ws.Range("D21").Value = "'2.009.54"
$ws.Range("E21").Value = "'  +3.27%  "

$ws.Range("D22").Value = "'4.574"
$ws.Range("E22").Value = "'  +3.03%  "

$ws.Range("D23").Value = "'8.690"
$ws.Range("E23").Value = "'  +1.49%  "

$ws.Range("D24").Value = "'5.232"
$ws.Range("E24").Value = "'  +2.56%  "

$ws.Range("D25").Value = "'141.49"
$ws.Range("E25").Value = "'  +2.40%  "

$ws.Range("D26").Value = "'15.72"
$ws.Range("E26").Value = "'  +2.33%  "

$ws.Range("D27").Value = "'1.859"
$ws.Range("E27").Value = "'  +6.62%  "

$ws.Range("D28").Value = "'109.23"
$ws.Range("E28").Value = "'  +2.88%  "

$ws.Range("D29").Value = "'1.402"
$ws.Range("E29").Value = "'  +0.09%  "

$ws.Range("D30").Value = "'4.201"
$ws.Range("E30").Value = "'  +7.23%  "

$ws.Range("D31").Value = "'0.08266"
$ws.Range("E31").Value = "'  +3.44%  "

$ws.Range("D32").Value = "'3.798"
$ws.Range("E32").Value = "'  +3.57%  "

$ws.Range("D33").Value = "'0.04856"
$ws.Range("E33").Value = "'  +7.92%  "

$ws.Range("D34").Value = "'1.070"
$ws.Range("E34").Value = "'  +6.76%  "

$ws.Range("D35").Value = "'2.611"
$ws.Range("E35").Value = "'  +0.17%  "

$ws.Range("D36").Value = "'0.6470"
$ws.Range("E36").Value = "'  +3.39%  "

$ws.Range("D37").Value = "'0.9427"
$ws.Range("E37").Value = "'  +0.56%  "

$ws.Range("D38").Value = "'2.576"
$ws.Range("E38").Value = "'  +6.28%  "

$ws.Range("D39").Value = "'2.043"
$ws.Range("E39").Value = "'  +2.07%  "

$ws.Range("D40").Value = "'5.946"
$ws.Range("E40").Value = "'  +6.67%  "

$ws.Range("D41").Value = "'0.01541"
$ws.Range("E41").Value = "'  +1.87%  "

$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "'  +0.40%  "

$ws.Range("D43").Value = "'99.98"
$ws.Range("E43").Value = "'  +0.48%  "

$ws.Range("D44").Value = "'0.3978"
$ws.Range("E44").Value = "'  +3.15%  "

$ws.Range("D45").Value = "'7.196"
$ws.Range("E45").Value = "'  +4.18%  "

$ws.Range("D47").Value = "'0.05414"
$ws.Range("E47").Value = "'  +0.62%  "

$ws.Range("D48").Value = "'8.016"
$ws.Range("E48").Value = "'  +2.76%  "

$ws.Range("E49").Value = "'  +4.80%  "

$ws.Range("D50").Value = "'30.62"
$ws.Range("E50").Value = "'  +0.91%  "

$ws.Range("D51").Value = "'52.88"
$ws.Range("E51").Value = "'  +2.33%  "
